# Update the "2024" sheet: a new September transaction entry was recorded,
# which shifts the existing September Details/Date rows (columns R:S) down
# by one row (rows 47-187), with the newest entry now occupying row 47.
# Separately, the "Broadband" label (column A) that lived on the last
# (blank) row of the sheet moves down to a brand new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# --- Shift September_Details (R) / September_Date (S) down by one row ---
# Walk from the bottom up so we never overwrite a value before reading it.
for ($row = 187; $row -ge 48; $row--) {
    $prevDetails = $ws.Cells.Item($row - 1, 18).Value()
    $prevDate = $ws.Cells.Item($row - 1, 19).Value()
    $ws.Cells.Item($row, 18).Value = $prevDetails
    $ws.Cells.Item($row, 19).Value = $prevDate
}

# New, most-recent entry goes into the now-vacated row 47
$ws.Cells.Item(47, 18).Value = "balance your axis"
$ws.Cells.Item(47, 19).Value = "2024-09-23 08:45:03"

# --- Move the "Broadband" label from row 195 down to a new row 196 ---
$srcRow = $ws.Range("A195:Y195")
$dstRow = $ws.Range("A196:Y196")
$srcRow.Copy($dstRow)
$ws.Cells.Item(195, 1).ClearContents()
